# matplotlib waterfall chart unit test
# Rename the measure-type labels in column C and tidy up the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "measure" column values to shorter labels.
# Order matters for shared-string table layout: introduce "abs" and "rel"
# before finally renaming the header cell to "Type".
$ws.Range("C2").Value = "abs"
$ws.Range("C3").Value = "rel"
$ws.Range("C4").Value = "rel"
$ws.Range("C5").Value = "rel"
$ws.Range("C7").Value = "rel"
$ws.Range("C8").Value = "rel"
$ws.Range("C9").Value = "rel"
$ws.Range("C1").Value = "Type"

# Remove the explicit (unused) style applied to the final total row label.
$ws.Range("A10").Style = "Normal"

# Update the view: select C2, and zoom out to 220%.
[void]$ws.Range("C2").Select()
$excel.ActiveWindow.Zoom = 220
